# Updates from WRI on 11/7 with data updates.
$wb = $excel.ActiveWorkbook

# Sheet 2 "BDSBaPCF" becomes the active/selected sheet (was sheet 1 "About").
$wsData = $wb.Worksheets.Item("BDSBaPCF")
$wsData.Activate()

# B17 ("natural gas peaker") was a formula `=B9`; it is replaced by the
# static value 0.
$wsData.Range("B17").Value = 0

# The selection on the BDSBaPCF sheet moves from E4 to B9.
$wsData.Range("B9").Select()
